$wb = $excel.ActiveWorkbook

# --- ev_charging_uc sheet: update the two comma-separated timeslice lists ---
$wsEv = $wb.Worksheets.Item("ev_charging_uc")
$wsEv.Range("C13").Value = "WaD,RaD,RaP,SaD,FaD,FaP,SaP,WaP"
$wsEv.Range("C14").Value = "FaP,SaP,RaP,WaP,SaN,WaN,FaN,RaN"

# --- re_profiles sheet: re-order the M4:O7 lookup block ---
$wsRe = $wb.Worksheets.Item("re_profiles")

$wsRe.Range("M4").Value = "S"
$wsRe.Range("N4").Value = 0.3412182463807702
$wsRe.Range("O4").Value = "hydro"

$wsRe.Range("M5").Value = "R"
$wsRe.Range("N5").Value = 0.40054629882545745
$wsRe.Range("O5").Value = "hydro"

$wsRe.Range("M6").Value = "W"
$wsRe.Range("N6").Value = 0.27238459437312207
$wsRe.Range("O6").Value = "hydro"

$wsRe.Range("M7").Value = "F"
$wsRe.Range("N7").Value = 0.18585086042065005
$wsRe.Range("O7").Value = "hydro"

$wb.Application.CalculateFull()
